$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.742.58'
$ws.Range('E2').Value = '  -1.61%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.570.11'
$ws.Range('E3').Value = '  -1.61%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '577.81'
$ws.Range('E5').Value = '  -2.44%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '188.73'
$ws.Range('E6').Value = '  -1.72%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.631'
$ws.Range('E7').Value = '  -3.27%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.565.23'
$ws.Range('E8').Value = '  -1.60%  '

$ws.Range('E9').Value = '  +0.10%  '

$ws.Range('E10').Value = '  -1.89%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.660'
$ws.Range('E11').Value = '  -0.75%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '55.60'
$ws.Range('E12').Value = '  -4.75%  '

$ws.Range('E13').Value = '  +1.70%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.60'
$ws.Range('E14').Value = '  -2.12%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.149.41'
$ws.Range('E15').Value = '  -1.33%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '19.80'
$ws.Range('E16').Value = '  +1.99%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.577.20'
$ws.Range('E17').Value = '  -1.44%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.751.82'
$ws.Range('E18').Value = '  -1.46%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.59'
$ws.Range('E19').Value = '  +0.13%  '

$ws.Range('E20').Value = '  -0.03%  '

$ws.Range('E21').Value = '  -1.43%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '473.58'
$ws.Range('E22').Value = '  -4.70%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '19.43'
$ws.Range('E23').Value = '  +12.14%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.03'
$ws.Range('E24').Value = '  -6.94%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '95.96'
$ws.Range('E25').Value = '  +5.31%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.37'
$ws.Range('E26').Value = '  -3.23%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.99'
$ws.Range('E27').Value = '  -4.39%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.01'
$ws.Range('E28').Value = '  -2.32%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.34'
$ws.Range('E29').Value = '  -1.55%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.76'
$ws.Range('E30').Value = '  +2.56%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '32.27'
$ws.Range('E31').Value = '  -0.52%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '12.19'
$ws.Range('E32').Value = '  -0.52%  '

$ws.Range('E33').Value = '  +0.29%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '66.09'
$ws.Range('E34').Value = '  +1.28%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '581.06'
$ws.Range('E35').Value = '  -6.79%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '38.88'
$ws.Range('E36').Value = '  +1.64%  '

$ws.Range('E37').Value = '  +0.01%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0₃0795'
$ws.Range('E38').Value = '  -4.46%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.395'
$ws.Range('E39').Value = '  -4.01%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.23'
$ws.Range('E40').Value = '  +18.85%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.89'
$ws.Range('E41').Value = '  +6.74%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.137'
$ws.Range('E42').Value = '  -6.68%  '

$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.45'
$ws.Range('E43').Value = '  -6.09%  '

$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.231.84'
$ws.Range('E44').Value = '  -2.87%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.08'
$ws.Range('E45').Value = '  +0.13%  '

$ws.Range('E46').Value = '  -1.39%  '

$ws.Range('E47').Value = '  +0.91%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.33'
$ws.Range('E48').Value = '  +1.58%  '

$ws.Range('E49').Value = '  -0.42%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.00'
$ws.Range('E50').Value = '  +0.03%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.12'
